$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tester B passed all three checks - fill in the result column (I) for rows 7-9
$ws.Range("I7").Value = "Pass"
$ws.Range("I8").Value = "Pass"
$ws.Range("I9").Value = "Pass"

# Move the active selection to I13 (matches the author's final cursor position)
$ws.Range("I13").Select()
